$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "Платные услуги (ДМС)" financing source now tagged as |pay|org|dms|
$ws.Range("F10").Value = "|pay|org|dms|"

# Row 12: "Платные физ. лица" financing source now tagged as |pay|individual|
$ws.Range("F12").Value = "|pay|individual|"

# Update the active selection to reflect the author's last click at F15
$ws.Range("F15").Select()
